$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.837.00"
$ws.Range("E2").Value = "  -0.12%  "

$ws.Range("D3").Value = "2.960.76"
$ws.Range("E3").Value = "  +0.87%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "106.78"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -4.72%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.557"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.64%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.610"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.49%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.17%  "

$ws.Range("E11").Value = "  +1.36%  "

$ws.Range("E12").Value = "  -4.49%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.05"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.12%  "

$ws.Range("D14").Value = "3.421.72"
$ws.Range("E14").Value = "  +0.91%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.58"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.50%  "

$ws.Range("D16").Value = "2.971.54"
$ws.Range("E16").Value = "  +1.54%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.988"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.17%  "

$ws.Range("D18").Value = "51.857.57"
$ws.Range("E18").Value = "  -0.23%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.40"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.59%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.44%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -5.31%  "

$ws.Range("E22").Value = "  -1.80%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "69.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.77%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.85"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.12%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.13%  "

$ws.Range("E26").Value = "  -3.87%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "26.64"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.25%  "

$ws.Range("E28").Value = "  +0.13%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.37"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.98%  "

$ws.Range("E30").Value = "  +1.77%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "10.28"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.77%  "

$ws.Range("E32").Value = "  -2.97%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "36.07"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "50.58"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.45%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0433"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.64%  "

$ws.Range("E37").Value = "  -0.03%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -4.52%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "17.78"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -5.09%  "

$ws.Range("E40").Value = "  -4.67%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.69"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("E42").Value = "  -1.39%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "123.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.17%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "22.42"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.85%  "

$ws.Range("E45").Value = "  -3.99%  "

$ws.Range("D46").Value = "2.112.22"
$ws.Range("E46").Value = "  -2.88%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.33"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.90%  "

$ws.Range("E48").Value = "  -8.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.236"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.23%  "

$ws.Range("E50").Value = "  -4.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.919"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -2.67%  "
